$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (A4 = 2)
$ws.Range("B4").Value = 11.47138019017781
$ws.Range("C4").Value = 16820.8
$ws.Range("D4").Value = 0.1148893923255633
$ws.Range("E4").Value = 165
$ws.Range("F4").Value = 373.4
$ws.Range("I4").Value = 0.2602785146347409
$ws.Range("J4").Value = 15.4
$ws.Range("K4").Value = 0.0102324902907245

# Row 5 (A5 = 3)
$ws.Range("B5").Value = 6.506973090568204
$ws.Range("C5").Value = 9534.4
$ws.Range("D5").Value = 0.04526183093699399
$ws.Range("E5").Value = 68.8
$ws.Range("F5").Value = 370.2
$ws.Range("I5").Value = 0.2520432498535652
$ws.Range("J5").Value = 12.8
$ws.Range("K5").Value = 0.008581153657621576
